{"js": "const body = context.document.body;\n\n// The contact-info text (\"email || phone || city, state\") was appended as its\n// own run right after the GitHub hyperlink. Find that exact run of text and\n// delete it so the website copy no longer exposes personal contact details.\nconst results = body.search(\" || qmeyer1995@gmail.com || 2604137437 || Saginaw, MI\", {\n  matchCase: true\n});\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# The contact-info text (\"email || phone || city, state\") lives in its own\n# run right after the GitHub hyperlink in the header block. Find that exact\n# text and replace it with nothing so the website copy no longer exposes\n# personal contact details.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \" || qmeyer1995@gmail.com || 2604137437 || Saginaw, MI\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"\"\n\n# wdReplaceAll = 2, wdFindContinue = 1\n$find.Execute(\n    $find.Text,        # FindText\n    $false,            # MatchCase\n    $false,            # MatchWholeWord\n    $false,            # MatchWildcards\n    $false,            # MatchSoundsLike\n    $false,            # MatchAllWordForms\n    $true,             # Forward\n    1,                 # Wrap (wdFindContinue)\n    $false,            # Format\n    $find.Replacement.Text,  # ReplaceWith\n    2                  # Replace (wdReplaceAll)\n)\n"}
